# The "2024" worksheet tracks monthly transaction details/timestamps.
# A new September transaction was recorded (most recent at the top of the
# September block), pushing every existing row at/after row 36 down by one
# row (dimension grows from A1:Y124 to A1:Y125).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at position 36, shifting rows 36:124 down to 37:125.
$ws.Rows("36:36").Insert()

# Populate the new row with the latest September transaction.
$ws.Range("R36").Value = "balance your axis"
$ws.Range("S36").Value = "2024-09-10 11:21:40"
